# Remove the two "Test"/"testtestsetstst" text boxes from slide 1,
# leaving only the background picture shape in place.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.Shapes.Item("TextBox 5").Delete()
$s.Shapes.Item("TextBox 7").Delete()
